# Word COM-interop script implementing the tracked change:
#  * new "High Level Usage and Overview" section (with "Pyro channels"
#    subsection) inserted after the "Easing the development..." bullet
#  * clean-up of stray <w:proofErr/> spell-check markers (merged back into
#    plain runs) in several places
#  * the <w:lastRenderedPageBreak/> marker moves from the start of the
#    "Inertial Measurement Unit" heading run into the middle of the V_BAT
#    paragraph, splitting that run
#  * a new trailing paragraph "Skhjahdgfkjhs" after the KiCAD heading

$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaRange([string]$needle) {
    # Locates the paragraph containing $needle and returns its full Range.
    $r = $d.Content
    $null = $r.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Expand(4) | Out-Null   # wdParagraph
    return $r
}

function Set-ParaInnerXml($paraRange, [string]$innerXml) {
    # Replaces everything inside a paragraph (keeping its own paragraph
    # mark / pPr / paraId / rsids intact) with the given <w:p>...</w:p> xml.
    $target = $d.Range($paraRange.Start, $paraRange.End - 1)
    $target.InsertXML("<w:p " + $ns.Substring(6) + ">" + $innerXml + "</w:p>")
}

function Insert-ParasAfter($paraRange, [string]$parasXml) {
    # Inserts one or more sibling <w:p>...</w:p> elements right after the
    # paragraph identified by $paraRange, without disturbing it.
    $point = $d.Range($paraRange.End - 1, $paraRange.End - 1)
    $point.InsertXML($parasXml)
}

# ---------------------------------------------------------------------
# 1. New "High Level Usage and Overview" + "Pyro channels" sections,
#    inserted right after "Easing the development process..." bullet.
# ---------------------------------------------------------------------
$anchor = Get-ParaRange("Easing the development process of avionics")

$newSections = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>High Level Usage and Overview</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">This section teaches members of other teams which </w:t></w:r><w:r><w:t>work</w:t></w:r><w:r><w:t xml:space="preserve"> with Hal-1 the functions of the board at a high level. The usage guide for the pyro channels, basic wiring from the board out, and selecting modes are covered.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Pyro channels</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">To use the pyro channels, electric matches (e-matches) are wired with one end into </w:t></w:r><w:r><w:t>the e-match’s respective OUT and IN</w:t></w:r><w:r><w:t xml:space="preserve"> screw terminal</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> on the board</w:t></w:r><w:r><w:t>. E-matches are non-directional, so it does not matter which end of the wire is screwed into the terminal. Both will produce the same results.</w:t></w:r></w:p>
"@

Insert-ParasAfter $anchor $newSections

# ---------------------------------------------------------------------
# 2. "128 KBytes flash memory" - drop the spell-check proofErr wrapper
#    around "KBytes" and merge the three runs into one.
# ---------------------------------------------------------------------
$p = Get-ParaRange("128 ")
Set-ParaInnerXml $p '<w:r><w:t>128 KBytes flash memory</w:t></w:r>'

# ---------------------------------------------------------------------
# 3. V_BAT paragraph: split the big run so <w:lastRenderedPageBreak/>
#    now sits right before "backup battery exists...".
# ---------------------------------------------------------------------
$p = Get-ParaRange("These three voltage pins are inputs")
$vbatInner = '<w:r><w:t xml:space="preserve">These three voltage pins are inputs. They accept the 3.3V in use for all modules on HAL-1. The difference between V_DD and V_DDA is that V_DDA is a separate voltage supply for the analog components within the chip itself. </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">The usage of the V_DDA pin is usually to allow a filtered input to be </w:t></w:r>' + `
    '<w:r><w:t>supplied but</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> is not implemented in HAL-1 in the interest of simplicity and quick functionality. </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">V_BAT is the backup power source, preferably from a 3.3 battery source, although the MCU can accept anything from 1.62V to 3.6V. If no </w:t></w:r>' + `
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">backup battery exists, as is the case for HAL-1, this pin may be connected to 3.3V power. </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Regardless of usage, V_DD, V_BAT, and V_DDA must </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">all </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">be connected to power. </w:t></w:r>'
Set-ParaInnerXml $p $vbatInner

# ---------------------------------------------------------------------
# 4. Remove the (now redundant) <w:lastRenderedPageBreak/> that used to
#    sit in front of "Inertial Measurement Unit (IMU): BNO085".
# ---------------------------------------------------------------------
$p = Get-ParaRange("Inertial Measurement Unit")
Set-ParaInnerXml $p '<w:r><w:t>Inertial Measurement Unit (IMU): BNO085</w:t></w:r>'

# ---------------------------------------------------------------------
# 5. "GPS: ublox MAX-M8Q" - drop proofErr wrapper, merge runs.
# ---------------------------------------------------------------------
$p = Get-ParaRange("GPS: ")
Set-ParaInnerXml $p '<w:r><w:t>GPS: ublox MAX-M8Q</w:t></w:r>'

# ---------------------------------------------------------------------
# 6. "Radio: Ebyte E22 900T22S" - drop proofErr wrapper, merge runs.
# ---------------------------------------------------------------------
$p = Get-ParaRange("Radio: ")
Set-ParaInnerXml $p '<w:r><w:t>Radio: Ebyte E22 900T22S</w:t></w:r>'

# ---------------------------------------------------------------------
# 7. "STM32 CubeMX and MCU Pinout" - drop proofErr wrapper, merge runs.
#    (search on the full phrase - "STM32 " alone also matches the
#    earlier "...based on STM32 system architecture" bullet)
# ---------------------------------------------------------------------
$p = Get-ParaRange("STM32 CubeMX")
Set-ParaInnerXml $p '<w:r><w:t>STM32 CubeMX and MCU Pinout</w:t></w:r>'

# ---------------------------------------------------------------------
# 8. "Shdjfg" - drop the proofErr wrapper around the lone run.
# ---------------------------------------------------------------------
$p = Get-ParaRange("Shdjfg")
Set-ParaInnerXml $p '<w:r><w:t>Shdjfg</w:t></w:r>'

# ---------------------------------------------------------------------
# 9. "KiCAD and Hardware Project Organization" - drop proofErr wrapper,
#    merge runs, and add the new trailing "Skhjahdgfkjhs" paragraph.
# ---------------------------------------------------------------------
$p = Get-ParaRange("KiCAD")
Set-ParaInnerXml $p '<w:r><w:t>KiCAD and Hardware Project Organization</w:t></w:r>'

$p = Get-ParaRange("KiCAD and Hardware Project Organization")
Insert-ParasAfter $p '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Skhjahdgfkjhs</w:t></w:r></w:p>'
